# Jogos_da_Semana_FlashScore_2024-11-28.xlsx : refresh odds snapshot and
# drop the South Korea K League 1 (Asan vs Daegu) fixture that is no longer
# part of this week's list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed odds values (row-by-row, matching the source diff)
# Row 2 - Santa Cruz vs Independiente (Bolivia)
$ws.Range("I2").Value = 3.1
$ws.Range("L2").Value = 3.75
$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 10
$ws.Range("O2").Value = 1.33
$ws.Range("P2").Value = 3.25
$ws.Range("Q2").Value = 2.1
$ws.Range("R2").Value = 1.7
$ws.Range("S2").Value = 1.44
$ws.Range("T2").Value = 2.63
$ws.Range("U2").Value = 1.8
$ws.Range("V2").Value = 1.91
$ws.Range("W2").Value = 7.5
$ws.Range("AC2").Value = 9.5
$ws.Range("AG2").Value = 251
$ws.Range("AJ2").Value = 12
$ws.Range("AL2").Value = 26
$ws.Range("AT2").Value = 2.63
$ws.Range("BB2").Value = 201

# Row 3 - Kerala Blasters vs Goa (India)
$ws.Range("Q3").Value = 1.67
$ws.Range("R3").Value = 2.15

# Row 5 - Tacuary vs Sportivo Trinidense (Paraguay)
$ws.Range("G5").Value = 3.9
$ws.Range("H5").Value = 3.75
$ws.Range("I5").Value = 1.8
$ws.Range("M5").Value = 1.04
$ws.Range("N5").Value = 13
$ws.Range("Q5").Value = 1.83
$ws.Range("R5").Value = 1.98
$ws.Range("U5").Value = 1.8
$ws.Range("V5").Value = 1.91
$ws.Range("W5").Value = 12
$ws.Range("X5").Value = 21
$ws.Range("AA5").Value = 34
$ws.Range("AC5").Value = 12
$ws.Range("AD5").Value = 7.5
$ws.Range("AH5").Value = 7.5
$ws.Range("AI5").Value = 9
$ws.Range("AK5").Value = 15

# Row 6 - Al Feiha vs Al Orubah (Saudi Arabia)
$ws.Range("Q6").Value = 2.2
$ws.Range("R6").Value = 1.65

# Row 7
$ws.Range("G7").Value = 1.83
$ws.Range("I7").Value = 3.9
$ws.Range("M7").Value = 1.03
$ws.Range("N7").Value = 9.5
$ws.Range("Q7").Value = 2
$ws.Range("R7").Value = 1.8
$ws.Range("W7").Value = 7
$ws.Range("X7").Value = 8.5
$ws.Range("AB7").Value = 29
$ws.Range("AC7").Value = 9.5
$ws.Range("AE7").Value = 17
$ws.Range("AJ7").Value = 15
$ws.Range("AN7").Value = 3.75
$ws.Range("AU7").Value = 8.5
$ws.Range("AX7").Value = 23
$ws.Range("AZ7").Value = 81

# Row 8 (final row, was followed by the row to be removed)
$ws.Range("G8").Value = 2.35
$ws.Range("H8").Value = 3.25
$ws.Range("I8").Value = 2.8
$ws.Range("M8").Value = 1.05
$ws.Range("N8").Value = 8.5
$ws.Range("Q8").Value = 1.93
$ws.Range("R8").Value = 1.88
$ws.Range("Z8").Value = 23
$ws.Range("AM8").Value = 29
$ws.Range("AX8").Value = 15
$ws.Range("AY8").Value = 23

# Remove the last data row (row 9, South Korea - K League 1 / Asan vs Daegu)
# entirely; this shifts nothing below it, it's simply dropped.
$ws.Rows.Item(9).Delete()
